# Generate Report for Handoff
# Updates the "ht" (handoff) priority marker and the latest handoff / handback
# generation timestamps for the files that were just re-handed-off.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsOverview = $wb.Worksheets.Item("Overview")

# zh-cn sheet: rows 6 (34febcc8...) and 7 (9a56f03f...) were just handed off again.
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("H6").Value = "2017-02-28 06:34:35"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("H7").Value = "2017-02-28 06:34:35"

# de-de sheet: rows 4 (24fdbda1...), 5 (34febcc8...), 6 (9a56f03f...) and 7 (ae402bb3...)
# are marked with the same "ht" priority.
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("E7").Value = "ht"

# Overview sheet: "Latest HO Xliff Generate Date" column for the affected files
# moves forward to reflect the newly generated handoff report.
$wsOverview.Range("G4").Value = "2017-02-28 06:34:51"
$wsOverview.Range("G5").Value = "2017-02-28 06:34:51"
$wsOverview.Range("G6").Value = "2017-02-28 06:34:51"
$wsOverview.Range("G7").Value = "2017-02-28 06:34:51"
